# Apply corrected area/cost/count assumptions to the "Dwellings_buildings"
# mapping-scheme sheet (non-residential commercial occupancy classes).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dwellings_buildings")

# Row -> classification (col B), average_area (C), average_unit_cost (D),
# dwellings_per_building (F). classification_proportion (E) is unchanged.
$rows = @(
    @{ Row = 2;  B = "Wholesale and retail trade";          C = 130;  D = 318.7451140065147; F = 1 },
    @{ Row = 3;  B = "Professional and technical services";  C = 130;  D = 366.1182410423453; F = 1 },
    @{ Row = 4;  B = "All other services";                   C = 130;  D = 330.4421824104235; F = 1 },
    @{ Row = 5;  B = "Wholesale and retail trade";          C = 260;  D = 318.7451140065147; F = 2 },
    @{ Row = 6;  B = "Professional and technical services";  C = 260;  D = 366.1182410423453; F = 2 },
    @{ Row = 7;  B = "All other services";                   C = 260;  D = 330.4421824104235; F = 2 },
    @{ Row = 8;  B = "Wholesale and retail trade";          C = 450;  D = 318.7451140065147; F = 3 },
    @{ Row = 9;  B = "Professional and technical services";  C = 450;  D = 366.1182410423453; F = 3 },
    @{ Row = 10; B = "All other services";                   C = 450;  D = 330.4421824104235; F = 3 },
    @{ Row = 11; B = "Wholesale and retail trade";          C = 900;  D = 318.7451140065147; F = 5 },
    @{ Row = 12; B = "Professional and technical services";  C = 900;  D = 366.1182410423453; F = 5 },
    @{ Row = 13; B = "All other services";                   C = 900;  D = 330.4421824104235; F = 5 },
    @{ Row = 14; B = "Professional and technical services";  C = 1200; D = 479.5798045602606; F = 5 },
    @{ Row = 15; B = "Professional and technical services";  C = 1200; D = 318.7451140065147; F = 5 },
    @{ Row = 16; B = "All other services";                   C = 1200; D = 330.4421824104235; F = 5 },
    @{ Row = 17; B = "Professional and technical services";  C = 3200; D = 479.5798045602606; F = 10 },
    @{ Row = 18; B = "Professional and technical services";  C = 3200; D = 318.7451140065147; F = 10 },
    @{ Row = 19; B = "All other services";                   C = 3200; D = 330.4421824104235; F = 10 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}
